# HoaDon1.xlsx edit: the "HoaDon2" detail-row placeholders are replaced with
# "HoaDon1" placeholders (also re-ordering STT/TenPhong/GiaPhong/SoNgay/TongTien
# into their natural column order B..F), and the active selection moves from
# G15 to F15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 holds the merge-field placeholders for the detail line of the
# invoice table. Previously they referenced "HoaDon2"; now they reference
# "HoaDon1" (and are written back in plain B->F column order).
$ws.Range("B15").Value = "%HoaDon1.STT;insert:copystyles"
$ws.Range("C15").Value = "%HoaDon1.TenPhong"
$ws.Range("D15").Value = "%HoaDon1.GiaPhong"
$ws.Range("E15").Value = "%HoaDon1.SoNgay"
$ws.Range("F15").Value = "%HoaDon1.TongTien"

# Move the active selection from G15 to F15.
$ws.Range("F15").Select()
